$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.002991666666666667
$ws.Range("H2").Value = 0.008975
$ws.Range("I2").Value = 0.0003566413595017623
$ws.Range("J2").Value = 0.0003566413595017623
$ws.Range("M2").Value = 15.24491733333333
$ws.Range("N2").Value = 45.73475199999999
$ws.Range("O2").Value = 0.4831257321597052
$ws.Range("P2").Value = 0.4831257321597052
$ws.Range("Q2").Value = 0.04560771102222222
$ws.Range("R2").Value = 0.4104693992
$ws.Range("S2").Value = 0.0001723026179277216
$ws.Range("T2").Value = 0.0001723026179277215
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.002991666666666667
$ws.Range("H3").Value = 0.008975
$ws.Range("I3").Value = 0.0003566413595017623
$ws.Range("J3").Value = 0.0003566413595017623
$ws.Range("O3").Value = 0.327710667227878
$ws.Range("P3").Value = 0.327710667227878
$ws.Range("Q3").Value = 0.0309363224
$ws.Range("R3").Value = 0.2784269016
$ws.Range("S3").Value = 0.00011687517788338
$ws.Range("T3").Value = 0.00011687517788338
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.002991666666666667
$ws.Range("H4").Value = 0.008975
$ws.Range("I4").Value = 0.0003566413595017623
$ws.Range("J4").Value = 0.0003566413595017623
$ws.Range("M4").Value = 5.969012333333333
$ws.Range("N4").Value = 17.907037
$ws.Range("O4").Value = 0.1891636006124168
$ws.Range("P4").Value = 0.1891636006124168
$ws.Range("Q4").Value = 0.01785729523055556
$ws.Range("R4").Value = 0.160715657075
$ws.Range("S4").Value = 0.00006746356369066073
$ws.Range("T4").Value = 0.00006746356369066073
$ws.Range("I5").Value = 0.9971069332391614
$ws.Range("J5").Value = 0.9971069332391616
$ws.Range("M5").Value = 15.24491733333333
$ws.Range("N5").Value = 45.73475199999999
$ws.Range("O5").Value = 0.4831257321597052
$ws.Range("P5").Value = 0.4831257321597052
$ws.Range("Q5").Value = 127.5111920080071
$ws.Range("R5").Value = 1147.600728072064
$ws.Range("S5").Value = 0.4817280171626882
$ws.Range("T5").Value = 0.4817280171626882
$ws.Range("I6").Value = 0.9971069332391614
$ws.Range("J6").Value = 0.9971069332391616
$ws.Range("O6").Value = 0.327710667227878
$ws.Range("P6").Value = 0.327710667227878
$ws.Range("Q6").Value = 86.492552622208
$ws.Range("R6").Value = 778.432973599872
$ws.Range("S6").Value = 0.3267625783893488
$ws.Range("T6").Value = 0.3267625783893489
$ws.Range("I7").Value = 0.9971069332391614
$ws.Range("J7").Value = 0.9971069332391616
$ws.Range("M7").Value = 5.969012333333333
$ws.Range("N7").Value = 17.907037
$ws.Range("O7").Value = 0.1891636006124168
$ws.Range("P7").Value = 0.1891636006124168
$ws.Range("Q7").Value = 49.92587766085376
$ws.Range("R7").Value = 449.3328989476839
$ws.Range("S7").Value = 0.1886163376871245
$ws.Range("T7").Value = 0.1886163376871245
$ws.Range("G8").Value = 0.02127666666666667
$ws.Range("H8").Value = 0.06383
$ws.Range("I8").Value = 0.002536425401336767
$ws.Range("J8").Value = 0.002536425401336767
$ws.Range("M8").Value = 15.24491733333333
$ws.Range("N8").Value = 45.73475199999999
$ws.Range("O8").Value = 0.4831257321597052
$ws.Range("P8").Value = 0.4831257321597052
$ws.Range("Q8").Value = 0.3243610244622222
$ws.Range("R8").Value = 2.919249220159999
$ws.Range("S8").Value = 0.0012254123790893
$ws.Range("T8").Value = 0.0012254123790893
$ws.Range("G9").Value = 0.02127666666666667
$ws.Range("H9").Value = 0.06383
$ws.Range("I9").Value = 0.002536425401336767
$ws.Range("J9").Value = 0.002536425401336767
$ws.Range("O9").Value = 0.327710667227878
$ws.Range("P9").Value = 0.327710667227878
$ws.Range("Q9").Value = 0.22001843552
$ws.Range("R9").Value = 1.98016591968
$ws.Range("S9").Value = 0.0008312136606458103
$ws.Range("T9").Value = 0.0008312136606458103
$ws.Range("G10").Value = 0.02127666666666667
$ws.Range("H10").Value = 0.06383
$ws.Range("I10").Value = 0.002536425401336767
$ws.Range("J10").Value = 0.002536425401336767
$ws.Range("M10").Value = 5.969012333333333
$ws.Range("N10").Value = 17.907037
$ws.Range("O10").Value = 0.1891636006124168
$ws.Range("P10").Value = 0.1891636006124168
$ws.Range("Q10").Value = 0.1270006857455555
$ws.Range("R10").Value = 1.14300617171
$ws.Range("S10").Value = 0.0004797993616016572
$ws.Range("T10").Value = 0.0004797993616016573
